$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (row 2 through 360) holds a "Förändrad" date value that was
# bulk-updated from serial 45186 (2023-09-17) to serial 45188 (2023-09-19)
# for every data row in the sheet.
$newValue = 45188
$lastRow = 360

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
